$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2736
$ws.Range("E2").Value = 226
$ws.Range("F2").Value = 226
$ws.Range("G2").Value = 280
$ws.Range("H2").Value = 215
$ws.Range("I2").Value = 147
$ws.Range("J2").Value = 69
$ws.Range("K2").Value = 5339
$ws.Range("L2").Value = 1571
$ws.Range("M2").Value = 3769
$ws.Range("N2").Value = 2959
$ws.Range("O2").Value = 810
$ws.Range("P2").Value = 279
$ws.Range("Q2").Value = 232
$ws.Range("R2").Value = -140
$ws.Range("S2").Value = -115
$ws.Range("T2").Value = 68
$ws.Range("U2").Value = 164
$ws.Range("V2").Value = 237
$ws.Range("W2").Value = 8.25
$ws.Range("X2").Value = 7.87
$ws.Range("Y2").Value = 5.02
$ws.Range("Z2").Value = 4.05
$ws.Range("AA2").Value = 41.68
$ws.Range("AB2").Value = 972.82
$ws.Range("AC2").Value = 263
$ws.Range("AD2").Value = 15.65
$ws.Range("AE2").Value = 5400
$ws.Range("AF2").Value = 0.76
$ws.Range("AG2").Value = 165
$ws.Range("AH2").Value = 4.01
$ws.Range("AI2").Value = 61.59
$ws.Range("AJ2").Value = 55895292

# Row 3
$ws.Range("D3").Value = 2609
$ws.Range("E3").Value = 234
$ws.Range("F3").Value = 234
$ws.Range("G3").Value = 105
$ws.Range("H3").Value = 55
$ws.Range("I3").Value = 85
$ws.Range("J3").Value = -31
$ws.Range("K3").Value = 5366
$ws.Range("L3").Value = 1684
$ws.Range("M3").Value = 3683
$ws.Range("N3").Value = 2910
$ws.Range("O3").Value = 773
$ws.Range("P3").Value = 279
$ws.Range("Q3").Value = 170
$ws.Range("R3").Value = -137
$ws.Range("S3").Value = 37
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = 93
$ws.Range("V3").Value = 349
$ws.Range("W3").Value = 8.99
$ws.Range("X3").Value = 2.1
$ws.Range("Y3").Value = 2.91
$ws.Range("Z3").Value = 1.02
$ws.Range("AA3").Value = 45.72
$ws.Range("AB3").Value = 966.69
$ws.Range("AC3").Value = 153
$ws.Range("AD3").Value = 24.21
$ws.Range("AE3").Value = 5310
$ws.Range("AF3").Value = 0.7
$ws.Range("AG3").Value = 140
$ws.Range("AH3").Value = 3.78
$ws.Range("AI3").Value = 89.82
$ws.Range("AJ3").Value = 55895292

# Row 4
$ws.Range("D4").Value = 2597
$ws.Range("E4").Value = 224
$ws.Range("F4").Value = 224
$ws.Range("G4").Value = 426
$ws.Range("H4").Value = 348
$ws.Range("I4").Value = 165
$ws.Range("J4").Value = 183
$ws.Range("K4").Value = 5543
$ws.Range("L4").Value = 1605
$ws.Range("M4").Value = 3938
$ws.Range("N4").Value = 3015
$ws.Range("O4").Value = 923
$ws.Range("P4").Value = 279
$ws.Range("Q4").Value = 198
$ws.Range("R4").Value = -263
$ws.Range("S4").Value = 35
$ws.Range("T4").Value = 226
$ws.Range("U4").Value = -28
$ws.Range("V4").Value = 367
$ws.Range("W4").Value = 8.63
$ws.Range("X4").Value = 13.39
$ws.Range("Y4").Value = 5.57
$ws.Range("Z4").Value = 6.38
$ws.Range("AA4").Value = 40.76
$ws.Range("AB4").Value = 998.08
$ws.Range("AC4").Value = 295
$ws.Range("AD4").Value = 11.34
$ws.Range("AE4").Value = 5502
$ws.Range("AF4").Value = 0.61
$ws.Range("AG4").Value = 140
$ws.Range("AH4").Value = 4.18
$ws.Range("AI4").Value = 46.46
$ws.Range("AJ4").Value = 55895292

# Row 5
$ws.Range("D5").Value = 2652
$ws.Range("E5").Value = 157
$ws.Range("F5").Value = 157
$ws.Range("G5").Value = 219
$ws.Range("H5").Value = 174
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 74
$ws.Range("K5").Value = 5527
$ws.Range("L5").Value = 1524
$ws.Range("M5").Value = 4003
$ws.Range("N5").Value = 3053
$ws.Range("O5").Value = 950
$ws.Range("P5").Value = 279
$ws.Range("Q5").Value = 297
$ws.Range("R5").Value = -121
$ws.Range("S5").Value = -139
$ws.Range("T5").Value = 73
$ws.Range("U5").Value = 224
$ws.Range("V5").Value = 349
$ws.Range("W5").Value = 5.92
$ws.Range("X5").Value = 6.57
$ws.Range("Y5").Value = 3.3
$ws.Range("Z5").Value = 3.15
$ws.Range("AA5").Value = 38.07
$ws.Range("AB5").Value = 1007.44
$ws.Range("AC5").Value = 179
$ws.Range("AD5").Value = 16.76
$ws.Range("AE5").Value = 5572
$ws.Range("AF5").Value = 0.54
$ws.Range("AG5").Value = 140
$ws.Range("AH5").Value = 4.67
$ws.Range("AI5").Value = 76.69
$ws.Range("AJ5").Value = 55895292

# Row 6
$ws.Range("D6").Value = 2418
$ws.Range("E6").Value = 93
$ws.Range("F6").Value = 93
$ws.Range("G6").Value = 77
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 11
$ws.Range("K6").Value = 4962
$ws.Range("L6").Value = 1371
$ws.Range("M6").Value = 3591
$ws.Range("N6").Value = 2999
$ws.Range("P6").Value = 279
$ws.Range("Q6").Value = 113
$ws.Range("R6").Value = -702
$ws.Range("S6").Value = 546
$ws.Range("T6").Value = 118
$ws.Range("U6").Value = -4
$ws.Range("V6").Value = 559
$ws.Range("W6").Value = 3.85
$ws.Range("X6").Value = 1.92
$ws.Range("Y6").Value = 0.36
$ws.Range("Z6").Value = 0.89
$ws.Range("AA6").Value = 38.17
$ws.Range("AB6").Value = 992.31
$ws.Range("AC6").Value = 19
$ws.Range("AD6").Value = 138.65
$ws.Range("AE6").Value = 5473
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 5.57
$ws.Range("AI6").Value = 756.52
$ws.Range("AJ6").Value = 55895292

# Row 7 - clear cells
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8 - clear cells
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9 - clear cells
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
